$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Amr Al Memari"
$summary.Range("B4").Value = 5095.97
$summary.Range("B6").Value = 574760
$summary.Range("B7").Value = 439863
$summary.Range("B8").Value = 134897
$summary.Range("B9").Value = 1.31

# --- Assets sheet ---
$assets = $wb.Worksheets.Item("Assets")
# insert a new row before current row 3 (Liquid Assets), shifting it + totals down
$assets.Range("A3").EntireRow.Insert()

# Update row 2 (Economy Car -> Mid-range Car, value)
$assets.Range("B2").Value = "Mid-range Car"
$assets.Range("C2").Value = 106823

# New row 3: Vehicles / Luxury Car / 461624 (copy formatting from row 2)
$assets.Range("A2:C2").Copy()
$assets.Range("A3:C3").PasteSpecial(-4122)
$assets.Range("A3").Value = "Vehicles"
$assets.Range("B3").Value = "Luxury Car"
$assets.Range("C3").Value = 461624

# Row 4 (previously row 3, Liquid Assets / Savings Account) - update value
$assets.Range("C4").Value = 6313

# Row 5 (previously row 4, TOTAL ASSETS)
$assets.Range("C5").Value = 574760

# --- Liabilities sheet ---
$liab = $wb.Worksheets.Item("Liabilities")
# insert two new rows before current row 3 (Credit Cards), shifting it + totals down
$liab.Range("A3:A4").EntireRow.Insert()

# Update row 2 values
$liab.Range("C2").Value = 64094
$liab.Range("D2").Value = 1335
$liab.Range("E2").Value = 4

# New row 3: Auto Loans / Vehicle Loan 2 / 276974 / 3297 / 7 (copy formatting from row 2)
$liab.Range("A2:E2").Copy()
$liab.Range("A3:E3").PasteSpecial(-4122)
$liab.Range("A3").Value = "Auto Loans"
$liab.Range("B3").Value = "Vehicle Loan 2"
$liab.Range("C3").Value = 276974
$liab.Range("D3").Value = 3297
$liab.Range("E3").Value = 7

# New row 4: Personal Loans / Personal Loan / 62156 / 1727 / 3 (copy formatting from row 2)
$liab.Range("A2:E2").Copy()
$liab.Range("A4:E4").PasteSpecial(-4122)
$liab.Range("A4").Value = "Personal Loans"
$liab.Range("B4").Value = "Personal Loan"
$liab.Range("C4").Value = 62156
$liab.Range("D4").Value = 1727
$liab.Range("E4").Value = 3

# Row 5 (previously row 3, Credit Cards) - update values
$liab.Range("C5").Value = 36639
$liab.Range("D5").Value = 1832

# Row 6 (previously row 4, TOTAL LIABILITIES)
$liab.Range("C6").Value = 439863
